$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Data")

# Insert a new first row and mark it with "[" (shifts existing data down by one row)
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "["

# Append the closing marker "]" right after the shifted data (now ends at row 66)
$ws.Range("A66").Value = "]"

# Leave the selection on the next empty row, as if the user just pressed Enter
$ws.Range("A67").Select()

# Make Test_Data the active sheet/tab
$ws.Activate()
